$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Switch the Runmode flag for CustomerSuite (B3) from "Y" to "N"
$ws.Range("B3").Value = "N"

# Move the active selection to B3 (matches the new active cell in the diff)
$ws.Range("B3").Select()
